$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 21000
$ws.Range("I7").Value = 7500
$ws.Range("J7").Value = 30000
$ws.Range("K7").Value = 7500
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = -7388
$ws.Range("N7").Value = -30224
$ws.Range("H14").Value = 21000
$ws.Range("I14").Value = 7500
$ws.Range("J14").Value = 30000
$ws.Range("K14").Value = 7500
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = -7309
$ws.Range("N14").Value = -30382
$ws.Range("H42").Value = 225.71428
$ws.Range("I42").Value = 70
$ws.Range("J42").Value = 433.33334
$ws.Range("K42").Value = 210
$ws.Range("L42").Value = 1300.00002
$ws.Range("M42").Value = 20
$ws.Range("N42").Value = -1760.00002
$ws.Range("H76").Value = 3047.5715
$ws.Range("I76").Value = 3011.7058
$ws.Range("K76").Value = 3011.7058
$ws.Range("M76").Value = -2696.7058
$ws.Range("H79").Value = 3047.5715
$ws.Range("I79").Value = 3011.7058
$ws.Range("K79").Value = 3011.7058
$ws.Range("M79").Value = -1919.7058
$ws.Range("H98").Value = 1333.0646
$ws.Range("I98").Value = 769.05554
$ws.Range("K98").Value = 769.05554
$ws.Range("M98").Value = 728.94446
$ws.Range("H121").Value = 2745.111
$ws.Range("J121").Value = 3050.625
$ws.Range("L121").Value = 9151.875
$ws.Range("N121").Value = -12645.875
$ws.Range("H122").Value = 1333.0646
$ws.Range("I122").Value = 769.05554
$ws.Range("K122").Value = 2307.16662
$ws.Range("M122").Value = 142.83338
$ws.Range("H137").Value = 1684.6111
$ws.Range("I137").Value = 1242.3572
$ws.Range("J137").Value = 3232.5
$ws.Range("K137").Value = 3727.0716
$ws.Range("L137").Value = 9697.5
$ws.Range("M137").Value = -1177.0716
$ws.Range("N137").Value = -14797.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H117").Value = 29983.166
$ws.Range("J117").Value = 29983.166
$ws.Range("L117").Value = 29983.166
$ws.Range("N117").Value = -39161.166
$ws.Range("H119").Value = 30033.334
$ws.Range("J119").Value = 30033.334
$ws.Range("L119").Value = 30033.334
$ws.Range("N119").Value = -39709.334
$ws.Range("H121").Value = 29997.777
$ws.Range("J121").Value = 29997.777
$ws.Range("L121").Value = 29997.777
$ws.Range("N121").Value = -33491.777

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 29893
$ws.Range("J112").Value = 29893
$ws.Range("L112").Value = 29893
$ws.Range("N112").Value = -32847
$ws.Range("H128").Value = 3000
$ws.Range("I128").Value = 3000
$ws.Range("K128").Value = 9000
$ws.Range("M128").Value = -6510

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 743
$ws.Range("I12").Value = 743
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 743
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -573
$ws.Range("N12").ClearContents()
$ws.Range("H23").Value = 26666.666
$ws.Range("J23").Value = 26666.666
$ws.Range("L23").Value = 26666.666
$ws.Range("N23").Value = -27146.666
$ws.Range("H27").Value = 26666.666
$ws.Range("J27").Value = 26666.666
$ws.Range("L27").Value = 26666.666
$ws.Range("N27").Value = -27050.666
$ws.Range("H31").Value = 19609510
$ws.Range("I31").Value = 37037960
$ws.Range("J31").Value = 2505.875
$ws.Range("K31").Value = 37037960
$ws.Range("L31").Value = 2505.875
$ws.Range("M31").Value = -37037665
$ws.Range("N31").Value = -3095.875
$ws.Range("H34").Value = 19609510
$ws.Range("I34").Value = 37037960
$ws.Range("J34").Value = 2505.875
$ws.Range("K34").Value = 37037960
$ws.Range("L34").Value = 2505.875
$ws.Range("M34").Value = -37037758
$ws.Range("N34").Value = -2909.875
$ws.Range("H134").Value = 25583104
$ws.Range("I134").Value = 2565850.2
$ws.Range("J134").Value = 250001330
$ws.Range("K134").Value = 7697550.600000001
$ws.Range("L134").Value = 750003990
$ws.Range("M134").Value = -7695015.600000001
$ws.Range("N134").Value = -750009060

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 794.5
$ws.Range("I17").Value = 194
$ws.Range("J17").Value = 2596
$ws.Range("K17").Value = 582
$ws.Range("L17").Value = 7788
$ws.Range("M17").Value = -413
$ws.Range("N17").Value = -8126
$ws.Range("H100").Value = 4050
$ws.Range("J100").Value = 4050
$ws.Range("L100").Value = 12150
$ws.Range("N100").Value = -13772
$ws.Range("H122").Value = 919.3570999999999
$ws.Range("I122").Value = 547.25
$ws.Range("J122").Value = 1849.625
$ws.Range("K122").Value = 4925.25
$ws.Range("L122").Value = 16646.625
$ws.Range("M122").Value = -2475.25
$ws.Range("N122").Value = -21546.625
$ws.Range("H131").Value = 903.36
$ws.Range("I131").Value = 822.25
$ws.Range("J131").Value = 906.73956
$ws.Range("K131").Value = 2466.75
$ws.Range("L131").Value = 2720.21868
$ws.Range("M131").Value = 2573.25
$ws.Range("N131").Value = -12800.21868
$ws.Range("H132").Value = 5389378
$ws.Range("I132").Value = 2383182.2
$ws.Range("J132").Value = 22224072
$ws.Range("K132").Value = 21448639.8
$ws.Range("L132").Value = 200016648
$ws.Range("M132").Value = -21446109.8
$ws.Range("N132").Value = -200021708

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 3401.4
$ws.Range("I17").Value = 504
$ws.Range("J17").Value = 5333
$ws.Range("K17").Value = 504
$ws.Range("L17").Value = 5333
$ws.Range("M17").Value = -336
$ws.Range("N17").Value = -5669
$ws.Range("H111").Value = 19096.5
$ws.Range("J111").Value = 19096.5
$ws.Range("L111").Value = 19096.5
$ws.Range("N111").Value = -25230.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 12505.167
$ws.Range("I19").Value = 261.75
$ws.Range("J19").Value = 36992
$ws.Range("K19").Value = 261.75
$ws.Range("L19").Value = 36992
$ws.Range("M19").Value = -91.75
$ws.Range("N19").Value = -37332
$ws.Range("H110").Value = 21849.625
$ws.Range("J110").Value = 21849.625
$ws.Range("L110").Value = 21849.625
$ws.Range("N110").Value = -30029.625
$ws.Range("H122").Value = 6990.4
$ws.Range("I122").Value = 9067.333000000001
$ws.Range("J122").Value = 3875
$ws.Range("K122").Value = 27201.999
$ws.Range("L122").Value = 11625
$ws.Range("M122").Value = -24751.999
$ws.Range("N122").Value = -16525

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32246
$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 150000
$ws.Range("N89").Value = -161232
$ws.Range("H119").Value = 30828.215
$ws.Range("J119").Value = 30828.215
$ws.Range("L119").Value = 30828.215
$ws.Range("N119").Value = -40504.215

Write-Host "Applied Belias_Profits price/profit data refresh across 8 sheets."